$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new columns at K:L (everything from the old K column onward
# shifts two columns to the right, e.g. old K -> new M, old O -> new Q).
$null = $ws.Range("K1:L1").EntireColumn.Insert()

# Widen the two new columns to fit the "Khoảng cách Check-in/out" headers.
$ws.Columns.Item(11).ColumnWidth = 23.666666666666664
$ws.Columns.Item(12).ColumnWidth = 24.333333333333336

# New "distance" header columns (row 7 = human-readable group header,
# row 8 = report-engine template placeholders), mirroring the style
# already used by their row thanks to the column insert above.
$ws.Range("K7").Value = "Khoảng cách Check-in"
$ws.Range("L7").Value = "Khoảng cách Check-out"
$ws.Range("K8").Value = "{{ReportStoreCheckeds.CheckInDistance}}"
$ws.Range("L8").Value = "{{ReportStoreCheckeds.CheckOutDistance}}"

$null = $ws.Range("G8").Select()
